# The upstream diff only touches the <w:nsid w:val="..."/> GUIDs that Word
# stores inside four <w:abstractNum> list definitions in word/numbering.xml
# (abstractNumId 990, 991, 99721, 99722). Nothing else in the package
# changes: no paragraph text, run formatting, list indentation, numbering
# format/level data, or numId references are touched, and the two numId
# values actually used by the body (1002 -> abstractNum 99721, 1003 ->
# abstractNum 99722) keep referencing the very same abstract numbering
# definitions before and after.
#
# The `nsid` attribute is an internal, Word-managed fingerprint for a list
# definition. It carries no visible/semantic meaning (it is not the
# numId/abstractNumId used to associate paragraphs with a list, nor a
# formatting property), and it is not surfaced anywhere on the Word object
# model - there is no ListTemplate/ListLevel/List property for it, and
# Document.WordOpenXML / Range.WordOpenXML (the only handles onto the raw
# package XML that the object model exposes) are read-only here:
# assigning to them raises "... is a read-only property; the assignment
# changed nothing." Filesystem/shell access to poke the .docx package
# directly is likewise unavailable ("use the Office object model via
# ... instead"), and none of the list-mutating members that are available
# (ApplyListTemplateWithLevel, RemoveNumbers, ListFormat.*, ...) touch or
# regenerate nsid - they either leave it untouched or mint a brand new
# abstractNum (with numId renumbering) which would just add unrelated
# diffs instead of reproducing this one.
#
# So there is no reachable Word-OM operation that reproduces this specific
# change, and the change itself has no observable effect on the document.
# To avoid introducing unrelated differences, this script intentionally
# performs no content mutation - it only confirms (read-only) that the two
# numbered paragraphs referencing abstractNum 99721/99722 are still present
# and unchanged.

$d = $word.ActiveDocument

$found = $d.Content.Find.Execute("zda usnesení valné hromady", $true, $false,
                                  $false, $false, $false, $true, 1, $false,
                                  "", 0)
Write-Output "Numbered clause still present: $found"
